$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.182.23'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '2.060.99'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '248.74'
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '57.17'
$ws.Range("E8").Value = '  -2.25%  '
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").Value = '0.0785'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '16.27'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '0.913'
$ws.Range("E13").Value = '  +13.37%  '
$ws.Range("D14").Value = '2.357.00'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").Value = '5.79'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").Value = '2.059.10'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '18.67'
$ws.Range("E17").Value = '  +12.71%  '
$ws.Range("D18").Value = '37.179.38'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").Value = '74.92'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").Value = '0.0₃0902'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").Value = '237.93'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '2.49'
$ws.Range("E24").Value = '  +4.27%  '
$ws.Range("D25").Value = '9.68'
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("D26").Value = '2.19'
$ws.Range("E26").Value = '  -4.05%  '
$ws.Range("D27").Value = '169.99'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").Value = '20.28'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").Value = '5.17'
$ws.Range("E30").Value = '  +8.70%  '
$ws.Range("D31").Value = '1.16'
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").Value = '4.64'
$ws.Range("E33").Value = '  +3.39%  '
$ws.Range("D34").Value = '0.0890'
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '2.28'
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("E37").Value = '  +0.73%  '
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("E39").Value = '  +17.09%  '
$ws.Range("D40").Value = '3.11'
$ws.Range("E40").Value = '  +8.65%  '
$ws.Range("D41").Value = '0.102'
$ws.Range("E41").Value = '  -11.64%  '
$ws.Range("D42").Value = '17.73'
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("D45").Value = '96.76'
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").Value = '1.277.09'
$ws.Range("E47").Value = '  -1.14%  '
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("D49").Value = '6.85'
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").Value = '2.250.89'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = '44.44'
$ws.Range("E51").Value = '  +1.11%  '
